$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @{Row=2; A=0; B=94; Img="dog/dog100.png"; Word="danken"; Cat="dog"},
    @{Row=3; A=1; B=42; Img="car/car092.png"; Word="dienen"; Cat="car"},
    @{Row=4; A=2; B=3; Img="dog/dog091.png"; Word="meinen"; Cat="dog"},
    @{Row=5; A=3; B=89; Img="car/car099.png"; Word="duschen"; Cat="car"},
    @{Row=6; A=4; B=49; Img="dog/dog089.png"; Word="stoßen"; Cat="dog"},
    @{Row=7; A=5; B=56; Img="dog/dog067.png"; Word="wachsen"; Cat="dog"},
    @{Row=8; A=6; B=115; Img="dog/dog087.png"; Word="grenzen"; Cat="dog"},
    @{Row=9; A=7; B=44; Img="dog/dog104.png"; Word="parken"; Cat="dog"},
    @{Row=10; A=8; B=126; Img="car/car124.png"; Word="holen"; Cat="car"},
    @{Row=11; A=9; B=32; Img="car/car093.png"; Word="atmen"; Cat="car"},
    @{Row=12; A=10; B=102; Img="dog/dog077.png"; Word="lernen"; Cat="dog"},
    @{Row=13; A=11; B=17; Img="car/car077.png"; Word="zögern"; Cat="car"},
    @{Row=14; A=12; B=86; Img="car/car115.png"; Word="süßen"; Cat="car"},
    @{Row=15; A=13; B=7; Img="dog/dog075.png"; Word="sparen"; Cat="dog"},
    @{Row=16; A=14; B=15; Img="dog/dog066.png"; Word="passen"; Cat="dog"},
    @{Row=17; A=15; B=30; Img="car/car064.png"; Word="bauen"; Cat="car"},
    @{Row=18; A=16; B=38; Img="dog/dog078.png"; Word="angeln"; Cat="dog"},
    @{Row=19; A=17; B=110; Img="dog/dog109.png"; Word="hassen"; Cat="dog"},
    @{Row=20; A=18; B=103; Img="car/car102.png"; Word="quellen"; Cat="car"},
    @{Row=21; A=19; B=13; Img="car/car084.png"; Word="kriegen"; Cat="car"},
    @{Row=22; A=20; B=85; Img="car/car095.png"; Word="rechnen"; Cat="car"},
    @{Row=23; A=21; B=43; Img="car/car082.png"; Word="heißen"; Cat="car"},
    @{Row=24; A=22; B=10; Img="car/car074.png"; Word="wählen"; Cat="car"},
    @{Row=25; A=23; B=50; Img="dog/dog064.png"; Word="wehen"; Cat="dog"},
    @{Row=26; A=24; B=45; Img="dog/dog088.png"; Word="wecken"; Cat="dog"},
    @{Row=27; A=25; B=39; Img="dog/dog090.png"; Word="tollen"; Cat="dog"},
    @{Row=28; A=26; B=76; Img="dog/dog112.png"; Word="lächeln"; Cat="dog"},
    @{Row=29; A=27; B=4; Img="car/car069.png"; Word="spüren"; Cat="car"},
    @{Row=30; A=28; B=107; Img="car/car096.png"; Word="hacken"; Cat="car"},
    @{Row=31; A=29; B=27; Img="dog/dog072.png"; Word="ärgern"; Cat="dog"},
    @{Row=32; A=30; B=25; Img="car/car094.png"; Word="münzen"; Cat="car"},
    @{Row=33; A=31; B=67; Img="car/car078.png"; Word="streifen"; Cat="car"}
)

foreach ($item in $data) {
    $r = $item.Row
    $ws.Cells.Item($r, 1).Value = $item.A
    $ws.Cells.Item($r, 2).Value = $item.B
    $ws.Cells.Item($r, 3).Value = $item.Img
    $ws.Cells.Item($r, 4).Value = $item.Word
    $ws.Cells.Item($r, 5).Value = $item.Cat
}
